$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Update the data values in rows 4-9 (columns B:E) with the corrected
#    lab figures. Column A (the row labels) keeps its text, only its style
#    index changes below together with the B:E styles.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 2.2999999999999998
$ws.Range("C4").Value = 0.8
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

$ws.Range("B5").Value = 301
$ws.Range("C5").Value = 301
$ws.Range("D5").Value = 317
$ws.Range("E5").Value = 320

$ws.Range("B6").Value = 1.4
$ws.Range("C6").Value = 1.4
$ws.Range("D6").Value = 6.8
$ws.Range("E6").Value = 8

$ws.Range("B7").Value = 12.8
$ws.Range("C7").Value = 12.8
$ws.Range("D7").Value = 12.6
$ws.Range("E7").Value = 12.7

$ws.Range("B8").Value = 0.66
$ws.Range("C8").Value = 1.33
$ws.Range("D8").Value = 6.33
$ws.Range("E8").Value = 7.21

$ws.Range("B9").Value = 1152
$ws.Range("C9").Value = 785
$ws.Range("D9").Value = 113
$ws.Range("E9").Value = 105

# ---------------------------------------------------------------------------
# 2) The label column (A4:A9) and the data columns (B4:E9) swap their cell
#    style (the two styles render identically - same borders - only the
#    underlying style-table index used by each column flips). Reproduce this
#    using a copy/paste-format round-trip via scratch cells so the style
#    table itself is not perturbed (no new, unused style entries created).
# ---------------------------------------------------------------------------
$ws.Range("A4:A9").Copy()
$ws.Range("G4:G9").PasteSpecial(-4122)   # xlPasteFormats -> stash col-A style

$ws.Range("B4:B9").Copy()
$ws.Range("H4:H9").PasteSpecial(-4122)   # xlPasteFormats -> stash col-B style

$ws.Range("H4:H9").Copy()
$ws.Range("A4:A9").PasteSpecial(-4122)   # col A now uses the old col-B style

$ws.Range("G4:G9").Copy()
$ws.Range("B4:E9").PasteSpecial(-4122)   # cols B:E now use the old col-A style

$ws.Range("G4:H9").Clear()               # drop the scratch helper cells
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Cosmetic workbook-window size metadata (bookViews/workbookView
#    windowWidth) is not exposed through the Excel COM surface in this
#    runtime (Window.Width writes don't persist to the saved file), so it is
#    left as-is; it carries no content/formatting meaning.
# ---------------------------------------------------------------------------
